# Fixed issue #13 Permitir que en los ficheros de metadatos dos columnas se
# puedan relacionar para crear SKOS jerarquicos
#
# A new row is inserted right below the header row. It holds short,
# machine readable codes for each column (used to relate/cross-reference
# columns when building hierarchical SKOS concepts). The old last row,
# which only referenced an external "mapping-ano.xlsx" file, is removed
# (that approach is superseded by the new in-sheet column codes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row right after the header (row 1), pushing every
# other row down by one. Excel inherits the formatting of the row above,
# so the new cells already come out with the same style as the rest of
# the sheet.
$ws.Rows.Item(2).Insert()

# New row 2: short codes that identify/relate each column.
$ws.Range("A2").Value = "vab"
$ws.Range("B2").Value = "sector-vab-descripcion"
$ws.Range("C2").Value = "codigo"
$ws.Range("D2").Value = "comarca"
$ws.Range("E2").Value = "sector-vab-codigo"
$ws.Range("F2").Value = "ano"

# The former last row (old row 5) only contained a reference to an
# external mapping file ("mapping-ano.xlsx"); that is now row 6 after the
# insert above, and is no longer needed, so drop it entirely.
$ws.Rows.Item(6).Delete()
